$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed error values in column B
$ws.Range("B8").Value = 0.4432709353245859
$ws.Range("B9").Value = -1.853424299602153
$ws.Range("B13").Value = 0.5217354724088806
$ws.Range("B16").Value = 0.308458428
$ws.Range("B18").Value = -0.1004302561888024
$ws.Range("B19").Value = 1.226733570319939
$ws.Range("B20").Value = 0.3320304208968201
$ws.Range("B21").Value = 0.6150051765378737
$ws.Range("B22").Value = -0.2683559768566441
$ws.Range("B23").Value = 0.2315426864241067

# Add two new rows for the new quarterly periods
$ws.Range("A24").Value = "2025-07-01_diff"
$ws.Range("A25").Value = "2025-10-01_diff"

# Copy the style (border/alignment/font) from an existing date cell (A23) to the new ones
$ws.Range("A23").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0
